$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 244; existing rows 244-327 shift down to 245-328
$ws.Rows.Item(244).Insert()

# Populate the newly inserted row 244 with its data
$ws.Cells.Item(244, 1).Value = 3
$ws.Cells.Item(244, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(244, 3).Value = "Coquimbo"
$ws.Cells.Item(244, 4).Value = 44627
$ws.Cells.Item(244, 5).Value = 5
$ws.Cells.Item(244, 6).Value = 100112031
$ws.Cells.Item(244, 7).Value = "Poroto verde"
$ws.Cells.Item(244, 8).Value = "Magnum"
$ws.Cells.Item(244, 9).Value = "Primera"
$ws.Cells.Item(244, 10).Value = 108
$ws.Cells.Item(244, 11).Value = 24000
$ws.Cells.Item(244, 12).Value = 25000
$ws.Cells.Item(244, 13).Value = 24454
$ws.Cells.Item(244, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(244, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(244, 16).Value = 978
$ws.Cells.Item(244, 17).Value = 25
$ws.Cells.Item(244, 18).Value = "Hortaliza"
